# CRM and pH calibration run 0404 dmb
# Adds a new data row (row 30) to the CRMAccuracyData sheet for the CRM
# bottle opened 2021-04-04, matching the existing table layout/formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = 20210404
$ws.Range("B30").Value = 2272.8870000000002
$ws.Range("C30").Value = 2231.4699999999998
$ws.Range("D30").Formula = "=100*(B30-C30)/C30"
$ws.Range("E30").Value = 180
$ws.Range("F30").Value = "CRM opened 20210404"

# Update the view to match where the user ended up after entering the
# new row: scrolled down a bit, with the next empty row selected.
$ws.Range("E31").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 26
$win.ScrollColumn = 1
